$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker list for column B, rows 2-49 (48 values), replacing the old list.
$bValues = @("NSE:ALANKIT","NSE:ASKAUTOLTD","NSE:AUROPHARMA","NSE:BALAJITELE","NSE:BHEL","NSE:COLPAL","NSE:DCW","NSE:DEEPAKFERT","NSE:DNAMEDIA","NSE:ERIS","NSE:ESTER","NSE:FILATEX","NSE:GMMPFAUDLR","NSE:HARIOMPIPE","NSE:IGARASHI","NSE:IMAGICAA","NSE:INOXGREEN","NSE:INOXWIND","NSE:J&KBANK","NSE:JINDALSAW","NSE:KANSAINER","NSE:KAYNES","NSE:KFINTECH","NSE:KHANDSE","NSE:LUXIND","NSE:MAKEINDIA","NSE:MANOMAY","NSE:MANUGRAPH","NSE:MAZDA","NSE:MOMOMENTUM","NSE:MOREPENLAB","NSE:NDLVENTURE","NSE:NEOGEN","NSE:NESCO","NSE:NETWORK18","NSE:NGLFINE","NSE:PANACEABIO","NSE:PARACABLES","NSE:PNB","NSE:PRECAM","NSE:PRSMJOHNSN","NSE:PYRAMID","NSE:RAIN","NSE:RECLTD","NSE:REPL","NSE:RHL","NSE:ROHLTD","NSE:RUPA")

# New ticker list for column C, rows 2-9 (8 values); rows 10-49 remain blank.
$cValues = @("NSE:AGRITECH","NSE:AXITA","NSE:COASTCORP","NSE:HERITGFOOD","NSE:INVENTURE","NSE:MURUDCERA","NSE:NESTLEIND","NSE:ORISSAMINE")

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}

# Rows 50-62 are removed entirely in the updated sheet (list now ends at row 49).
$ws.Range("A50:F62").EntireRow.Delete()

Write-Host "Edit complete"
